$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 2171.0417
$ws.Range("I15").Value = 2171.0417
$ws.Range("K15").Value = 6513.125100000001
$ws.Range("M15").Value = -6344.125100000001
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 1141.75
$ws.Range("I33").Value = 1760.5
$ws.Range("J33").Value = 523
$ws.Range("K33").Value = 1760.5
$ws.Range("L33").Value = 523
$ws.Range("M33").Value = -1531.5
$ws.Range("N33").Value = -981
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 4633.6665
$ws.Range("I43").Value = 4633.6665
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4633.6665
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4564.6665
$ws.Range("N43").ClearContents()
# Row 82: Rolling on Initiative
$ws.Range("H82").Value = 4571.75
$ws.Range("I82").Value = 3413.3333
$ws.Range("K82").Value = 10239.9999
$ws.Range("M82").Value = -9833.999899999999
# Row 85: Darkly Dreaming Dexterity (L)
$ws.Range("H85").Value = 4571.75
$ws.Range("I85").Value = 3413.3333
$ws.Range("K85").Value = 10239.9999
$ws.Range("M85").Value = -8835.999899999999
# Row 120: Supreme Official Strategy Guide
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1767.7693
$ws.Range("I132").Value = 1767.7693
$ws.Range("K132").Value = 5303.3079
$ws.Range("M132").Value = -2773.3079

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 12: Strait Ain't the Gate
$ws.Range("H12").Value = 15500
$ws.Range("I12").Value = 15500
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 15500
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -15327
$ws.Range("N12").ClearContents()
# Row 34: Insistent Sallets
$ws.Range("H34").Value = 7014
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 12028
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 12028
$ws.Range("M34").Value = -1729
$ws.Range("N34").Value = -12570
# Row 37: Get Shirty
$ws.Range("H37").Value = 30488
$ws.Range("J37").Value = 30488
$ws.Range("L37").Value = 30488
$ws.Range("N37").Value = -31034
# Row 41: Skillet Scandal
$ws.Range("H41").Value = 1236
$ws.Range("I41").Value = 883.2
$ws.Range("J41").Value = 3000
$ws.Range("K41").Value = 883.2
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = -469.2
$ws.Range("N41").Value = -3828
# Row 42: Kitty Get Your Helm
$ws.Range("H42").Value = 12031
$ws.Range("J42").Value = 12031
$ws.Range("L42").Value = 12031
$ws.Range("N42").Value = -13003
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 3750
$ws.Range("I45").Value = 4000
$ws.Range("K45").Value = 4000
$ws.Range("M45").Value = -3623
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 9373.955
$ws.Range("I74").Value = 10091.777
$ws.Range("K74").Value = 10091.777
$ws.Range("M74").Value = -9217.777
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 9373.955
$ws.Range("I77").Value = 10091.777
$ws.Range("K77").Value = 50458.885
$ws.Range("M77").Value = -46090.885
# Row 97: Ore for Me
$ws.Range("H97").Value = 1698.75
$ws.Range("I97").Value = 1698.75
$ws.Range("K97").Value = 1698.75
$ws.Range("M97").Value = -1202.75
# Row 104: See Shields by the Sea Shore
$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -51988
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 3160.3572
$ws.Range("I110").Value = 1480.3846
$ws.Range("K110").Value = 1480.3846
$ws.Range("M110").Value = 564.6153999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 46: Spice Cadet
$ws.Range("H46").Value = 6000
$ws.Range("J46").Value = 6000
$ws.Range("L46").Value = 6000
$ws.Range("N46").Value = -6596
# Row 94: High Steal
$ws.Range("H94").Value = 2040
$ws.Range("I94").Value = 2866.8333
$ws.Range("J94").Value = 799.75
$ws.Range("K94").Value = 2866.8333
$ws.Range("L94").Value = 799.75
$ws.Range("M94").Value = -2415.8333
$ws.Range("N94").Value = -1701.75
# Row 122: To Delight a Dancer
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 32: Daddy's Little Girl
$ws.Range("H32").Value = 1337.5
$ws.Range("I32").Value = 1337.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1337.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1021.5
$ws.Range("N32").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 191.4
$ws.Range("I2").Value = 103.5
$ws.Range("J2").Value = 323.25
$ws.Range("K2").Value = 621
$ws.Range("L2").Value = 1939.5
$ws.Range("M2").Value = -508
$ws.Range("N2").Value = -2165.5
# Row 10: A Real Fungi
$ws.Range("H10").Value = 417
$ws.Range("I10").Value = 678.3333
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 2034.9999
$ws.Range("L10").Value = 75
$ws.Range("M10").Value = -1895.9999
$ws.Range("N10").Value = -353
# Row 76: Old Victories, New Tastes
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
# Row 79: The Eats of Authenticity (L)
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
# Row 108: Meet for Meat
$ws.Range("H108").Value = 2077.5
$ws.Range("I108").Value = 436.33334
$ws.Range("K108").Value = 1309.00002
$ws.Range("M108").Value = 1570.99998
# Row 139: Najoothie
$ws.Range("H139").Value = 3029.5
$ws.Range("I139").Value = 3029.5
$ws.Range("K139").Value = 9088.5
$ws.Range("M139").Value = -3948.5
# Row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 3403.125
$ws.Range("I140").Value = 3460.7144
$ws.Range("K140").Value = 10382.1432
$ws.Range("M140").Value = -5202.143199999999
# Row 141: Ocean Explosion
$ws.Range("H141").Value = 3999
$ws.Range("I141").Value = 3999
$ws.Range("K141").Value = 11997
$ws.Range("M141").Value = -6817

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 29: Music to Their Ears
$ws.Range("H29").Value = 800
$ws.Range("J29").Value = 800
$ws.Range("L29").Value = 800
$ws.Range("N29").Value = -1380

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 11076.467
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590
# Row 27: Fire and Hide
$ws.Range("H27").Value = 11076.467
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 2512.5
$ws.Range("I46").Value = 2157.1428
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 2157.1428
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1969.1428
$ws.Range("N46").Value = -5376

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 33: I'll Be Your Wailer Today
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15500
# Row 36: Put a Lid on It
$ws.Range("H36").Value = 15000
$ws.Range("J36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15500
# Row 54: No Country for Cold Men
$ws.Range("H54").Value = 25443.875
$ws.Range("J54").Value = 25443.875
$ws.Range("L54").Value = 25443.875
$ws.Range("N54").Value = -26483.875
# Row 116: All-purpose Overgarments
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 3720.4443
$ws.Range("I136").Value = 3056.9412
$ws.Range("K136").Value = 9170.8236
$ws.Range("M136").Value = -6620.8236
